$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from I1 into the newly added J1 header cell (style s="1")
$ws.Range("I1").Copy($ws.Range("J1"))

$data = New-Object 'object[,]' 35,9
$data[0,0] = 0
$data[0,1] = 1
$data[0,2] = 2
$data[0,3] = 3
$data[0,4] = 4
$data[0,5] = 5
$data[0,6] = 6
$data[0,7] = 7
$data[0,8] = 8
$data[1,0] = 0
$data[1,1] = 0
$data[1,2] = 0
$data[1,3] = 0
$data[1,4] = 0
$data[1,5] = 0
$data[1,6] = 0.004614370468029005
$data[1,7] = 0.01840659340659341
$data[1,8] = 0
$data[2,0] = 0
$data[2,1] = 0.2462574850299392
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 0
$data[2,5] = 0.1555306007860756
$data[2,6] = 0.002636783124588002
$data[2,7] = 0.02499999999999998
$data[2,8] = 0
$data[3,0] = 0.05585831062670298
$data[3,1] = 0
$data[3,2] = 0.0915492957746479
$data[3,3] = 0.1081632653061223
$data[3,4] = 0.06666666666666667
$data[3,5] = 0
$data[3,6] = 0
$data[3,7] = 0
$data[3,8] = 0
$data[4,0] = 0
$data[4,1] = 0.01796407185628743
$data[4,2] = 0
$data[4,3] = 0
$data[4,4] = 0
$data[4,5] = 0.05277933745087039
$data[4,6] = 0.04878048780487806
$data[4,7] = 0.05302197802197791
$data[4,8] = 0
$data[5,0] = 0
$data[5,1] = 0
$data[5,2] = 0
$data[5,3] = 0
$data[5,4] = 0
$data[5,5] = 0
$data[5,6] = 0.003955174686882004
$data[5,7] = 0.02829670329670327
$data[5,8] = 0
$data[6,0] = 0.1539509536784742
$data[6,1] = 0
$data[6,2] = 0
$data[6,3] = 0
$data[6,4] = 0
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 0
$data[6,8] = 0.01456310679611651
$data[7,0] = 0
$data[7,1] = 0.2425149700598794
$data[7,2] = 0
$data[7,3] = 0
$data[7,4] = 0
$data[7,5] = 0.07804604154969133
$data[7,6] = 0.01582069874752802
$data[7,7] = 0.06071428571428557
$data[7,8] = 0
$data[8,0] = 0.01362397820163487
$data[8,1] = 0
$data[8,2] = 0
$data[8,3] = 0
$data[8,4] = 0
$data[8,5] = 0
$data[8,6] = 0.01450230718523402
$data[8,7] = 0
$data[8,8] = 0.2168284789644004
$data[9,0] = 0.1103542234332427
$data[9,1] = 0
$data[9,2] = 0.07746478873239437
$data[9,3] = 0.120408163265306
$data[9,4] = 0.06666666666666667
$data[9,5] = 0
$data[9,6] = 0
$data[9,7] = 0
$data[9,8] = 0
$data[10,0] = 0
$data[10,1] = 0
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[10,5] = 0
$data[10,6] = 0.01582069874752802
$data[10,7] = 0
$data[10,8] = 0.001618122977346278
$data[11,0] = 0.01771117166212533
$data[11,1] = 0
$data[11,2] = 0
$data[11,3] = 0
$data[11,4] = 0
$data[11,5] = 0
$data[11,6] = 0.01252471984179301
$data[11,7] = 0
$data[11,8] = 0.02831715210355989
$data[12,0] = 0
$data[12,1] = 0
$data[12,2] = 0
$data[12,3] = 0
$data[12,4] = 0
$data[12,5] = 0
$data[12,6] = 0.009887936717205011
$data[12,7] = 0.01236263736263737
$data[12,8] = 0
$data[13,0] = 0
$data[13,1] = 0
$data[13,2] = 0
$data[13,3] = 0
$data[13,4] = 0
$data[13,5] = 0.02189781021897811
$data[13,6] = 0.001977587343441002
$data[13,7] = 0.01263736263736265
$data[13,8] = 0
$data[14,0] = 0
$data[14,1] = 0
$data[14,2] = 0
$data[14,3] = 0
$data[14,4] = 0
$data[14,5] = 0
$data[14,6] = 0
$data[14,7] = 0
$data[14,8] = 0
$data[15,0] = 0.1294277929155316
$data[15,1] = 0
$data[15,2] = 0
$data[15,3] = 0.02244897959183674
$data[15,4] = 0
$data[15,5] = 0
$data[15,6] = 0.2096242584047469
$data[15,7] = 0.02252747252747252
$data[15,8] = 0.1844660194174753
$data[16,0] = 0
$data[16,1] = 0
$data[16,2] = 0
$data[16,3] = 0
$data[16,4] = 0
$data[16,5] = 0
$data[16,6] = 0
$data[16,7] = 0.02747252747252745
$data[16,8] = 0
$data[17,0] = 0.008174386920980926
$data[17,1] = 0
$data[17,2] = 0
$data[17,3] = 0
$data[17,4] = 0
$data[17,5] = 0
$data[17,6] = 0
$data[17,7] = 0
$data[17,8] = 0.009708737864077673
$data[18,0] = 0
$data[18,1] = 0
$data[18,2] = 0
$data[18,3] = 0
$data[18,4] = 0
$data[18,5] = 0
$data[18,6] = 0.027686222808174
$data[18,7] = 0.007967032967032965
$data[18,8] = 0
$data[19,0] = 0
$data[19,1] = 0
$data[19,2] = 0
$data[19,3] = 0
$data[19,4] = 0
$data[19,5] = 0
$data[19,6] = 0.007251153592617008
$data[19,7] = 0.1236263736263733
$data[19,8] = 0
$data[20,0] = 0
$data[20,1] = 0
$data[20,2] = 0
$data[20,3] = 0
$data[20,4] = 0
$data[20,5] = 0.003368893879842785
$data[20,6] = 0.0586684245220831
$data[20,7] = 0.1140109890109887
$data[20,8] = 0.001618122977346278
$data[21,0] = 0
$data[21,1] = 0
$data[21,2] = 0
$data[21,3] = 0
$data[21,4] = 0
$data[21,5] = 0
$data[21,6] = 0
$data[21,7] = 0
$data[21,8] = 0
$data[22,0] = 0
$data[22,1] = 0.1796407185628741
$data[22,2] = 0
$data[22,3] = 0
$data[22,4] = 0
$data[22,5] = 0.110612015721505
$data[22,6] = 0.02636783124588
$data[22,7] = 0.1057692307692305
$data[22,8] = 0
$data[23,0] = 0.04495912806539508
$data[23,1] = 0
$data[23,2] = 0
$data[23,3] = 0
$data[23,4] = 0
$data[23,5] = 0
$data[23,6] = 0
$data[23,7] = 0
$data[23,8] = 0.0008090614886731392
$data[24,0] = 0
$data[24,1] = 0
$data[24,2] = 0
$data[24,3] = 0
$data[24,4] = 0
$data[24,5] = 0
$data[24,6] = 0.04152933421226104
$data[24,7] = 0
$data[24,8] = 0.05016181229773459
$data[25,0] = 0
$data[25,1] = 0
$data[25,2] = 0
$data[25,3] = 0
$data[25,4] = 0
$data[25,5] = 0
$data[25,6] = 0
$data[25,7] = 0
$data[25,8] = 0.0008090614886731392
$data[26,0] = 0.0217983651226158
$data[26,1] = 0
$data[26,2] = 0
$data[26,3] = 0
$data[26,4] = 0
$data[26,5] = 0
$data[26,6] = 0.00922874093605801
$data[26,7] = 0
$data[26,8] = 0.01941747572815535
$data[27,0] = 0
$data[27,1] = 0
$data[27,2] = 0
$data[27,3] = 0
$data[27,4] = 0
$data[27,5] = 0
$data[27,6] = 0.01911667765326301
$data[27,7] = 0.02307692307692307
$data[27,8] = 0
$data[28,0] = 0
$data[28,1] = 0
$data[28,2] = 0
$data[28,3] = 0
$data[28,4] = 0
$data[28,5] = 0
$data[28,6] = 0.03823335530652602
$data[28,7] = 0.001648351648351648
$data[28,8] = 0
$data[29,0] = 0.001362397820163488
$data[29,1] = 0
$data[29,2] = 0
$data[29,3] = 0
$data[29,4] = 0
$data[29,5] = 0
$data[29,6] = 0.05207646671061308
$data[29,7] = 0
$data[29,8] = 0.02265372168284791
$data[30,0] = 0.0326975476839237
$data[30,1] = 0
$data[30,2] = 0
$data[30,3] = 0
$data[30,4] = 0
$data[30,5] = 0
$data[30,6] = 0
$data[30,7] = 0
$data[30,8] = 0.07200647249190935
$data[31,0] = 0.1444141689373299
$data[31,1] = 0
$data[31,2] = 0.8239436619718301
$data[31,3] = 0.5755102040816362
$data[31,4] = 0.8666666666666666
$data[31,5] = 0
$data[31,6] = 0
$data[31,7] = 0
$data[31,8] = 0.06310679611650477
$data[32,0] = 0.1934604904632151
$data[32,1] = 0
$data[32,2] = 0.007042253521126761
$data[32,3] = 0.1734693877551018
$data[32,4] = 0
$data[32,5] = 0
$data[32,6] = 0.003295978905735003
$data[32,7] = 0
$data[32,8] = 0.1221682847896442
$data[33,0] = 0
$data[33,1] = 0
$data[33,2] = 0
$data[33,3] = 0
$data[33,4] = 0
$data[33,5] = 0
$data[33,6] = 0.0006591957811470006
$data[33,7] = 0.01153846153846155
$data[33,8] = 0
$data[34,0] = 0.009536784741144413
$data[34,1] = 0
$data[34,2] = 0
$data[34,3] = 0
$data[34,4] = 0
$data[34,5] = 0
$data[34,6] = 0.001318391562294001
$data[34,7] = 0
$data[34,8] = 0.04288025889967636

$ws.Range("B1:J35").Value = $data
